$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2891.182
$ws.Range("I51").Value = 1834
$ws.Range("K51").Value = 1834
$ws.Range("M51").Value = -1350
$ws.Range("H86").Value = 4285.75
$ws.Range("I86").Value = 2552
$ws.Range("K86").Value = 2552
$ws.Range("M86").Value = -1429
$ws.Range("H89").Value = 4285.75
$ws.Range("I89").Value = 2552
$ws.Range("K89").Value = 12760
$ws.Range("M89").Value = -7144
$ws.Range("H107").Value = 1813.625
$ws.Range("I107").Value = 1415.1111
$ws.Range("K107").Value = 1415.1111
$ws.Range("M107").Value = 504.8888999999999
$ws.Range("H113").Value = 100001400
$ws.Range("I113").Value = 100001400
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 100001400
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -99998146
$ws.Range("N113").ClearContents()
$ws.Range("H125").Value = 1554.7142
$ws.Range("I125").Value = 2537
$ws.Range("J125").Value = 245
$ws.Range("K125").Value = 22833
$ws.Range("L125").Value = 2205
$ws.Range("M125").Value = -20373
$ws.Range("N125").Value = -7125
$ws.Range("H132").Value = 6175292.5
$ws.Range("I132").Value = 7754400.5
$ws.Range("K132").Value = 23263201.5
$ws.Range("M132").Value = -23260671.5
$ws.Range("H137").Value = 1607
$ws.Range("I137").Value = 1276
$ws.Range("J137").Value = 3427.5
$ws.Range("K137").Value = 3828
$ws.Range("L137").Value = 10282.5
$ws.Range("M137").Value = -1278
$ws.Range("N137").Value = -15382.5
$ws.Range("H138").Value = 1752.83
$ws.Range("I138").Value = 713.4375
$ws.Range("J138").Value = 1950.8096
$ws.Range("K138").Value = 2140.3125
$ws.Range("L138").Value = 5852.4288
$ws.Range("M138").Value = 2999.6875
$ws.Range("N138").Value = -16132.4288
$ws.Range("H141").Value = 1779.2307
$ws.Range("I141").Value = 1779.2307
$ws.Range("K141").Value = 5337.6921
$ws.Range("M141").Value = -157.6921000000002
# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2427.639
$ws.Range("I32").Value = 2543.5574
$ws.Range("K32").Value = 2543.5574
$ws.Range("M32").Value = -2256.5574
$ws.Range("H61").Value = 1224.6316
$ws.Range("I61").Value = 1027.5294
$ws.Range("K61").Value = 1027.5294
$ws.Range("M61").Value = -815.5293999999999
$ws.Range("H123").Value = 75000
$ws.Range("J123").Value = 75000
$ws.Range("L123").Value = 75000
$ws.Range("N123").Value = -84800
$ws.Range("H136").Value = 1224.6316
$ws.Range("I136").Value = 1027.5294
$ws.Range("K136").Value = 3082.5882
$ws.Range("M136").Value = -532.5881999999997
# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 24235.5
$ws.Range("I26").Value = 24235.5
$ws.Range("K26").Value = 24235.5
$ws.Range("M26").Value = -23943.5
# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 9940061
$ws.Range("I6").Value = 19870120
$ws.Range("J6").Value = 10002
$ws.Range("K6").Value = 19870120
$ws.Range("L6").Value = 10002
$ws.Range("M6").Value = -19870007
$ws.Range("N6").Value = -10228
$ws.Range("H31").Value = 972.1316
$ws.Range("I31").Value = 692.1818
$ws.Range("J31").Value = 2819.8
$ws.Range("K31").Value = 692.1818
$ws.Range("L31").Value = 2819.8
$ws.Range("M31").Value = -397.1818
$ws.Range("N31").Value = -3409.8
$ws.Range("H34").Value = 972.1316
$ws.Range("I34").Value = 692.1818
$ws.Range("J34").Value = 2819.8
$ws.Range("K34").Value = 692.1818
$ws.Range("L34").Value = 2819.8
$ws.Range("M34").Value = -490.1818
$ws.Range("N34").Value = -3223.8
$ws.Range("H44").Value = 600
$ws.Range("I44").Value = 600
$ws.Range("K44").Value = 600
$ws.Range("M44").Value = -158
$ws.Range("H58").Value = 750.9286
$ws.Range("I58").Value = 750
$ws.Range("J58").Value = 756.5
$ws.Range("K58").Value = 750
$ws.Range("L58").Value = 756.5
$ws.Range("M58").Value = -547
$ws.Range("N58").Value = -1162.5
$ws.Range("H132").Value = 10579.25
$ws.Range("I132").Value = 11795.2
$ws.Range("K132").Value = 35385.60000000001
$ws.Range("M132").Value = -32855.60000000001
$ws.Range("H134").Value = 11906179
$ws.Range("I134").Value = 15152765
$ws.Range("K134").Value = 45458295
$ws.Range("M134").Value = -45455760
$ws.Range("H135").Value = 53740
$ws.Range("J135").Value = 53740
$ws.Range("L135").Value = 53740
$ws.Range("N135").Value = -63880
$ws.Range("H136").Value = 750.9286
$ws.Range("I136").Value = 750
$ws.Range("J136").Value = 756.5
$ws.Range("K136").Value = 2250
$ws.Range("L136").Value = 2269.5
$ws.Range("M136").Value = 300
$ws.Range("N136").Value = -7369.5
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()
# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 139.66667
$ws.Range("I26").Value = 107.2
$ws.Range("J26").Value = 302
$ws.Range("K26").Value = 321.6
$ws.Range("L26").Value = 906
$ws.Range("M26").Value = -33.60000000000002
$ws.Range("N26").Value = -1482
$ws.Range("H32").Value = 1975
$ws.Range("H44").Value = 1800
$ws.Range("J44").Value = 3500
$ws.Range("L44").Value = 10500
$ws.Range("N44").Value = -11296
$ws.Range("H107").Value = 50776
$ws.Range("I107").Value = 1553
$ws.Range("J107").Value = 99999
$ws.Range("K107").Value = 4659
$ws.Range("L107").Value = 299997
$ws.Range("M107").Value = -2739
$ws.Range("N107").Value = -303837
$ws.Range("H131").Value = 11364822
$ws.Range("I131").Value = 142857820
$ws.Range("J131").Value = 1229.9753
$ws.Range("K131").Value = 428573460
$ws.Range("L131").Value = 3689.9259
$ws.Range("M131").Value = -428568420
$ws.Range("N131").Value = -13769.9259
$ws.Range("H132").Value = 647.5
$ws.Range("I132").Value = 647.5
$ws.Range("K132").Value = 5827.5
$ws.Range("M132").Value = -3297.5
# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 20235.5
$ws.Range("I99").Value = 20471
$ws.Range("J99").Value = 20000
$ws.Range("K99").Value = 20471
$ws.Range("L99").Value = 20000
$ws.Range("M99").Value = -18225
$ws.Range("N99").Value = -24492
$ws.Range("H122").Value = 4210.778
$ws.Range("I122").Value = 2632.3333
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 7896.999899999999
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -5446.999899999999
$ws.Range("N122").Value = -19900
$ws.Range("H133").Value = 41675.4
$ws.Range("J133").Value = 41675.4
$ws.Range("L133").Value = 41675.4
$ws.Range("N133").Value = -51795.4
# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2600
$ws.Range("I7").Value = 2600
$ws.Range("K7").Value = 2600
$ws.Range("M7").Value = -2488
$ws.Range("H22").Value = 987.625
$ws.Range("I22").Value = 551
$ws.Range("J22").Value = 1133.1666
$ws.Range("K22").Value = 551
$ws.Range("L22").Value = 1133.1666
$ws.Range("M22").Value = -256
$ws.Range("N22").Value = -1723.1666
$ws.Range("H27").Value = 987.625
$ws.Range("I27").Value = 551
$ws.Range("J27").Value = 1133.1666
$ws.Range("K27").Value = 551
$ws.Range("L27").Value = 1133.1666
$ws.Range("M27").Value = -444
$ws.Range("N27").Value = -1347.1666
$ws.Range("H46").Value = 5868.125
$ws.Range("I46").Value = 681.6667
$ws.Range("K46").Value = 681.6667
$ws.Range("M46").Value = -493.6667
$ws.Range("H55").Value = 839.9
$ws.Range("I55").Value = 300.33334
$ws.Range("J55").Value = 1071.1428
$ws.Range("K55").Value = 300.33334
$ws.Range("L55").Value = 1071.1428
$ws.Range("M55").Value = -127.33334
$ws.Range("N55").Value = -1417.1428
$ws.Range("H61").Value = 4999
$ws.Range("I61").Value = 4998
$ws.Range("K61").Value = 4998
$ws.Range("M61").Value = -4796
$ws.Range("H68").Value = 1589.2
$ws.Range("I68").Value = 1313.4286
$ws.Range("J68").Value = 2232.6667
$ws.Range("K68").Value = 1313.4286
$ws.Range("L68").Value = 2232.6667
$ws.Range("M68").Value = -564.4286
$ws.Range("N68").Value = -3730.6667
$ws.Range("H71").Value = 1589.2
$ws.Range("I71").Value = 1313.4286
$ws.Range("J71").Value = 2232.6667
$ws.Range("K71").Value = 6567.143
$ws.Range("L71").Value = 11163.3335
$ws.Range("M71").Value = -2823.143
$ws.Range("N71").Value = -18651.3335
$ws.Range("H113").Value = 4999
$ws.Range("I113").Value = 4998
$ws.Range("K113").Value = 4998
$ws.Range("M113").Value = -2828
$ws.Range("H126").Value = 2600
$ws.Range("I126").Value = 2600
$ws.Range("K126").Value = 7800
$ws.Range("M126").Value = -5330
$ws.Range("H136").Value = 13704.5
$ws.Range("I136").Value = 15486.571
$ws.Range("J136").Value = 1230
$ws.Range("K136").Value = 46459.713
$ws.Range("L136").Value = 3690
$ws.Range("M136").Value = -43909.713
$ws.Range("N136").Value = -8790
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()
# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H106").Value = 20377
$ws.Range("J106").Value = 20377
$ws.Range("L106").Value = 20377
$ws.Range("N106").Value = -22901
$ws.Range("H123").Value = 53752.332
$ws.Range("J123").Value = 53752.332
$ws.Range("L123").Value = 53752.332
$ws.Range("N123").Value = -63552.332
$ws.Range("H132").Value = 2059.186
$ws.Range("I132").Value = 1861.6765
$ws.Range("J132").Value = 2805.3333
$ws.Range("K132").Value = 5585.029500000001
$ws.Range("L132").Value = 8415.999899999999
$ws.Range("M132").Value = -3055.029500000001
$ws.Range("N132").Value = -13475.9999
$ws.Range("H136").Value = 587.8108
$ws.Range("I136").Value = 317.34784
$ws.Range("K136").Value = 952.0435200000001
$ws.Range("M136").Value = 1597.95648
